{"js": "// The author (1) finished the first answer's last sentence \u2014 which had\n// been split into two runs straddling Word's auto \"last edit\" bookmark\n// (_GoBack) \u2014 and (2) appended a closing remark to the final answer's\n// last sentence, which is where _GoBack now lives after the edit.\n\nconst body = context.document.body;\n\n// --- Edit 1: complete \"...\u5b9a\u4e49\u4e86\u5bc4\u5b58\" -> \"...\u5b9a\u4e49\u4e86\u5bc4\u5b58\u5668\uff0c\u5185\u5b58\u7b49\u7684\u6307\u4ee4\u8868\u793a\u3002\"\n// The sentence already reads complete in the paragraph's combined text\n// (the run split is invisible at the text level), so resolving the\n// search hit and writing the full sentence back merges the two runs.\nconst completedSentence = \"\u5b9a\u4e49\u4e86\u5bc4\u5b58\u5668\uff0c\u5185\u5b58\u7b49\u7684\u6307\u4ee4\u8868\u793a\u3002\";\nconst r1 = body.search(completedSentence, { matchCase: false, matchWholeWord: false });\nr1.load(\"text\");\nawait context.sync();\nif (r1.items.length > 0) {\n  r1.items[0].insertText(completedSentence, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Edit 2: append the new closing remark after \"...\u5b9e\u73b0\u4e86\u8fc7\u7a0b\u8c03\u7528\u3002\"\nconst addition = \"\uff08\u5373\u6808\u7684\u4f7f\u7528\u786e\u4fdd\u4e86\u5bc4\u5b58\u5668\u5171\u7528\u60c5\u51b5\u4e0b\uff0c\u4e0d\u4f1a\u56e0\u4e3a\u5207\u6362\u8fc7\u7a0b\u4e22\u5931\u72b6\u6001\uff09\";\nconst r2 = body.search(\"\u5b9e\u73b0\u4e86\u8fc7\u7a0b\u8c03\u7528\u3002\", { matchCase: false, matchWholeWord: false });\nr2.load(\"text\");\nawait context.sync();\n\nif (r2.items.length > 0) {\n  const target = r2.items[r2.items.length - 1];\n  const inserted = target.insertText(addition, Word.InsertLocation.end);\n  await context.sync();\n\n  // Relocate Word's hidden _GoBack bookmark (tracks the most recent edit\n  // location) so it again sits just before the trailing \"\uff09\", matching\n  // where Word itself leaves it after typing this text.\n  const closingParen = inserted.search(\"\uff09\", { matchCase: true });\n  closingParen.load(\"text\");\n  await context.sync();\n\n  if (closingParen.items.length > 0) {\n    const bmRange = closingParen.items[closingParen.items.length - 1].getRange(Word.RangeLocation.start);\n    await context.sync();\n\n    context.document.deleteBookmark(\"_GoBack\");\n    bmRange.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "# The author (1) finished the first answer's last sentence - which had\n# been split into two runs straddling Word's auto \"last edit\" bookmark\n# (_GoBack) - and (2) appended a closing remark to the final answer's\n# last sentence, which is where _GoBack now lives after the edit.\n\n$doc = $word.ActiveDocument\n\n# --- Edit 1: complete \"...\u5b9a\u4e49\u4e86\u5bc4\u5b58\" -> \"...\u5b9a\u4e49\u4e86\u5bc4\u5b58\u5668\uff0c\u5185\u5b58\u7b49\u7684\u6307\u4ee4\u8868\u793a\u3002\"\n# Re-typing the already-complete sentence over the Find hit (which spans\n# the old run split / bookmark) merges it back into a single run.\n$completedSentence = \"\u5b9a\u4e49\u4e86\u5bc4\u5b58\u5668\uff0c\u5185\u5b58\u7b49\u7684\u6307\u4ee4\u8868\u793a\u3002\"\n$rng1 = $doc.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Replacement.ClearFormatting()\n$rng1.Find.Execute($completedSentence, $false, $false, $false, $false, $false, $true, 1, $false, $completedSentence, 2) | Out-Null\n\n# --- Edit 2: append the new closing remark after \"...\u5b9e\u73b0\u4e86\u8fc7\u7a0b\u8c03\u7528\u3002\"\n$addition = \"\uff08\u5373\u6808\u7684\u4f7f\u7528\u786e\u4fdd\u4e86\u5bc4\u5b58\u5668\u5171\u7528\u60c5\u51b5\u4e0b\uff0c\u4e0d\u4f1a\u56e0\u4e3a\u5207\u6362\u8fc7\u7a0b\u4e22\u5931\u72b6\u6001\uff09\"\n$rng2 = $doc.Content\n$rng2.Find.ClearFormatting()\n$found = $rng2.Find.Execute(\"\u5b9e\u73b0\u4e86\u8fc7\u7a0b\u8c03\u7528\u3002\")\nif ($found) {\n    $rng2.Collapse(0)\n    $rng2.InsertAfter($addition)\n\n    # Relocate Word's hidden _GoBack bookmark (tracks the most recent\n    # edit location) so it again sits just before the trailing \"\uff09\",\n    # matching where Word itself leaves it after typing this text.\n    $bmPos = $rng2.End - 1\n    $bmRng = $doc.Range($bmPos, $bmPos)\n    $doc.Bookmarks.Add(\"_GoBack\", $bmRng) | Out-Null\n}\n"}
